$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("Q2").Value = 1.9
$ws.Range("R2").Value = 2
$ws.Range("G3").Value = 1.95
$ws.Range("Y3").Value = 10
$ws.Range("Z3").Value = 17
$ws.Range("AN3").Value = 3.75
$ws.Range("AO3").Value = 12
$ws.Range("AZ3").Value = 101
$ws.Range("M4").Value = 1.11
$ws.Range("N4").Value = 6.5
$ws.Range("G5").Value = 2.9
$ws.Range("I5").Value = 2.38
$ws.Range("J5").Value = 3.5
$ws.Range("M5").Value = 1.06
$ws.Range("N5").Value = 10
$ws.Range("U5").Value = 1.75
$ws.Range("V5").Value = 2
$ws.Range("W5").Value = 9
$ws.Range("AB5").Value = 29
$ws.Range("AK5").Value = 23
$ws.Range("AW5").Value = 4.5
$ws.Range("N6").Value = 13.9
$ws.Range("O6").Value = 1.12
$ws.Range("P6").Value = 4.45
$ws.Range("U6").Value = 1.83
$ws.Range("V6").Value = 1.93
$ws.Range("N7").Value = 7
$ws.Range("Q8").Value = 1.98
$ws.Range("R8").Value = 1.88
$ws.Range("O9").Value = 1.14
$ws.Range("P9").Value = 5.5
$ws.Range("G10").Value = 2.25
$ws.Range("N10").Value = 8
$ws.Range("O10").Value = 1.4
$ws.Range("P10").Value = 2.75
$ws.Range("Q10").Value = 2.3
$ws.Range("R10").Value = 1.6
$ws.Range("X10").Value = 10
$ws.Range("AF10").Value = 51
$ws.Range("AJ10").Value = 12
$ws.Range("AK10").Value = 34
$ws.Range("AN10").Value = 4.33
$ws.Range("AX10").Value = 19
$ws.Range("AY10").Value = 29
$ws.Range("BB10").Value = 251
$ws.Range("H11").Value = 3.1
$ws.Range("I11").Value = 3.9
$ws.Range("L11").Value = 4.75
$ws.Range("M11").Value = 1.11
$ws.Range("N11").Value = 6.5
$ws.Range("W11").Value = 5.5
$ws.Range("Z11").Value = 19
$ws.Range("AI11").Value = 17
$ws.Range("G13").Value = 1.27
$ws.Range("H13").Value = 5.25
$ws.Range("K13").Value = 2.6
$ws.Range("L13").Value = 9.5
$ws.Range("M13").Value = 1.03
$ws.Range("N13").Value = 15
$ws.Range("Q13").Value = 1.67
$ws.Range("R13").Value = 2.15
$ws.Range("S13").Value = 1.3
$ws.Range("T13").Value = 3.4
$ws.Range("W13").Value = 7
$ws.Range("AC13").Value = 12
$ws.Range("AD13").Value = 10
$ws.Range("AJ13").Value = 34
$ws.Range("AN13").Value = 3.1
$ws.Range("AT13").Value = 3.4
$ws.Range("AX13").Value = 51
$ws.Range("AY13").Value = 51
